# Generate Report for Archive
# The handoff file dbf4aa4f-437c-47fb-a18d-ea7794ef6d6f.md has been picked up,
# so the d5a60e5f-9ff4-49fb-9232-3e4af63a9d74.md row moves from
# "Ready for handoff" to "In Translation" on every sheet that tracks status.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status columns are "zh-cn" (E) and "de-de" (F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = "In Translation"
$wsOverview.Range("F6").Value = "In Translation"

# --- zh-cn sheet: Status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = "In Translation"

# --- de-de sheet: Status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = "In Translation"
